$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of environment data (order matters for shared-string indices)
$ws.Range("B3").Value = "https://i-preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do"
$ws.Range("A3").Value = "i-preproducciongestion.segurossura.com.ar"
$ws.Range("C3").Value = "su"
$ws.Range("D3").Value = "silverarrow"

# Add hyperlink on B3 pointing to itself (same URL as its text)
$ws.Hyperlinks.Add($ws.Range("B3"), "https://i-preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do")

# Update selection to match the post-edit state (C4 active cell)
$ws.Range("C4").Select()
